$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1341-1349 already exist with a blank (inlineStr) "backup" column R;
# fill them in with the now-known numeric backup flag (1 for the first row,
# 0 for the rest).
$backupValues = @(1, 0, 0, 0, 0, 0, 0, 0, 0)
$startRow = 1341
for ($i = 0; $i -lt $backupValues.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 18).Value = $backupValues[$i]
}

# Newly scraped trading days (2024-12-17 .. 2025-01-20), appended as rows 1350-1373.
# Column F (Adj Close) and column R (backup) are left blank for these rows, matching
# the source data which has not been back-filled yet.
$newRows = @(
    @(45643,2385,2399,2361.800048828125,2367.199951171875,138887,2024,12,17,0,0,0,51,0,0,0),
    @(45644,2367.949951171875,2389.64990234375,2346.10009765625,2359.89990234375,253905,2024,12,18,0,0,0,51,0,0,0),
    @(45645,2322,2370,2322,2363.39990234375,171380,2024,12,19,0,0,0,51,0,0,0),
    @(45646,2360,2373,2242.75,2249.949951171875,272969,2024,12,20,0,0,0,51,0,0,0),
    @(45649,2265,2275.5,2241.85009765625,2252.199951171875,161773,2024,12,23,0,0,0,52,0,0,0),
    @(45650,2259.949951171875,2275,2235.60009765625,2245.449951171875,143630,2024,12,24,0,0,0,52,0,0,0),
    @(45652,2250,2251.85009765625,2199.89990234375,2231.5,288582,2024,12,26,0,0,0,52,0,0,0),
    @(45653,2234,2271.75,2230,2256.35009765625,280034,2024,12,27,0,0,0,52,0,0,0),
    @(45656,2262,2282,2244,2268.800048828125,481181,2024,12,30,0,0,0,1,0,0,0),
    @(45657,2268.800048828125,2275,2226.14990234375,2248.800048828125,270289,2024,12,31,0,0,0,1,0,0,0),
    @(45658,2254.85009765625,2291.85009765625,2245.14990234375,2284.800048828125,75411,2025,1,1,0,0,0,1,0,0,0),
    @(45659,2288,2310.449951171875,2275.050048828125,2307.75,54293,2025,1,2,0,0,0,1,0,0,0),
    @(45660,2319,2338.050048828125,2244.449951171875,2250.39990234375,236058,2025,1,3,0,0,0,1,0,0,0),
    @(45663,2273.050048828125,2275,2165.550048828125,2176.449951171875,295756,2025,1,6,0,0,0,2,2,0,0),
    @(45664,2218,2312.699951171875,2188.89990234375,2293.550048828125,671612,2025,1,7,0,0,0,2,0,0,0),
    @(45665,2301,2316,2272.75,2299.64990234375,254717,2025,1,8,0,0,0,2,0,0,0),
    @(45666,2301,2371,2295.050048828125,2322.300048828125,239180,2025,1,9,0,0,0,2,0,0,0),
    @(45667,2315.050048828125,2339.35009765625,2279.10009765625,2325.14990234375,214208,2025,1,10,0,0,0,2,0,0,0),
    @(45670,2314.10009765625,2314.949951171875,2234.10009765625,2242.699951171875,144279,2025,1,13,0,0,0,3,0,0,0),
    @(45671,2245.35009765625,2290.39990234375,2228.25,2254.699951171875,350675,2025,1,14,0,0,0,3,0,0,0),
    @(45672,2266,2278.60009765625,2198,2209.64990234375,263290,2025,1,15,0,0,0,3,0,0,0),
    @(45673,2220.699951171875,2254.35009765625,2213.14990234375,2239.050048828125,221073,2025,1,16,0,0,0,3,0,0,0),
    @(45674,2245,2290,2224.699951171875,2282.60009765625,249982,2025,1,17,0,0,0,3,0,0,0),
    @(45677,2286,2311.449951171875,2239.10009765625,2270,210157,2025,1,20,0,0,0,4,0,0,0),
)

$colForIndex = @(1,2,3,4,5,7,8,9,10,11,12,13,14,15,16,17)  # A,B,C,D,E,G,H,I,J,K,L,M,N,O,P,Q (F and R stay blank)
$r = 1350
foreach ($row in $newRows) {
    for ($j = 0; $j -lt $colForIndex.Length; $j++) {
        $ws.Cells.Item($r, $colForIndex[$j]).Value = $row[$j]
    }
    # Column A carries the same "YYYY-MM-DD HH:MM:SS" datetime format as the rest of the sheet.
    $ws.Range("A$r").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $r++
}
